$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.618874073028564
$ws.Range("B1").Value = 2.667764186859131
$ws.Range("C1").Value = 3.001014471054077
$ws.Range("D1").Value = 3.359756708145142
$ws.Range("E1").Value = 2.340937852859497
